$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the contents of row 3 (the "ESC" row): A3:E3
$ws.Range("A3:E3").ClearContents()

# Move the active selection to A3 (matches the saved selection in the file)
$ws.Range("A3").Select()
